$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

$dateSerial = 46048
$startRow = 332

$data = @(
    @('ASAHI DEPOSITORY LLC Registered', 0, 0, 0, 0, 0, 0),
    @('ASAHI DEPOSITORY LLC Eligible', 0, 0, 0, 0, 0, 0),
    @('BRINK''S, INC. Registered', 87949.747, 0, 0, 0, 0, 87949.747),
    @('BRINK''S, INC. Eligible', 30578.352, 0, 0, 0, 0, 30578.352),
    @('CNT DEPOSITORY, INC. Registered', 1246.06, 0, 0, 0, 0, 1246.06),
    @('CNT DEPOSITORY, INC. Eligible', 0, 0, 0, 0, 0, 0),
    @('DELAWARE DEPOSITORY Registered', 1633.941, 0, 0, 0, 0, 1633.941),
    @('DELAWARE DEPOSITORY Eligible', 18459.584, 0, 0, 0, 0, 18459.584),
    @('HSBC BANK, USA Registered', 1394.758, 0, 0, 0, 0, 1394.758),
    @('HSBC BANK, USA Eligible', 9281.978999999999, 0, 0, 0, 0, 9281.978999999999),
    @('INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered', 2395.448, 0, 0, 0, 0, 2395.448),
    @('INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible', 0, 0, 0, 0, 0, 0),
    @('JP MORGAN CHASE BANK NA Registered', 114985.579, 0, 0, 0, 0, 114985.579),
    @('JP MORGAN CHASE BANK NA Eligible', 135413.823, 0, 0, 0, 0, 135413.823),
    @('LOOMIS INTERNATIONAL (US) LLC Registered', 63745.991, 0, 0, 0, 0, 63745.991),
    @('LOOMIS INTERNATIONAL (US) LLC Eligible', 132077.206, 0, 0, 0, 0, 132077.206),
    @('MALCA-AMIT USA, LLC Registered', 395.145, 0, 0, 0, 0, 395.145),
    @('MALCA-AMIT USA, LLC Eligible', 0, 0, 0, 0, 0, 0),
    @('MANFRA, TORDELLA & BROOKES, LLC Registered', 50220.42, 0, 0, 0, 0, 50220.42),
    @('MANFRA, TORDELLA & BROOKES, LLC Eligible', 1271.373, 0, 0, 0, 0, 1271.373),
    @('STONEX PRECIOUS METALS LLC Registered', 14122.765, 0, 0, 0, 0, 14122.765),
    @('STONEX PRECIOUS METALS LLC Eligible', 16.075, 0, 0, 0, 0, 16.075)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $dateSerial
    $aCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
    $ws.Cells.Item($row, 8).Value = $rec[6]
}
